# Insert a new price-report row at row 13 (shifting the existing rows
# 13..77 down to 14..78, and growing the used range to A1:R78), then
# populate the new row 13 with the reported values for this week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13..77 down to 14..78, opening up a blank row 13.
$ws.Rows(13).Insert()

# New data row (columns that stay constant for every "Camote" record on
# this sheet - Mercado ID/Mercado/Región/Codreg/Categoría ID/Categoría/
# Variedad/Calidad/Clasificación - are carried over explicitly as well).
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44677
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100114002
$ws.Cells.Item(13, 7).Value = "Camote"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 20
$ws.Cells.Item(13, 11).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 18000
$ws.Cells.Item(13, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 900
$ws.Cells.Item(13, 17).Value = 20
$ws.Cells.Item(13, 18).Value = "Hortaliza"
